# "Update countries & provincias Spain"
#
# The "Pais" sheet lists per-country COVID-19 stats (B=Casos totales,
# C=Nuevos casos, D=Casos activos, E=Recuperados, F=Casos criticos,
# G=Muertes hoy, H=Muertes), sorted by B descending. This refresh updates
# the timestamp and several countries' figures; three pairs/groups of
# rows swap rank (and therefore swap which country name sits on which
# row) because their updated totals reordered them relative to their
# neighbours.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp: 16:35 -> 17:05
$ws.Range("A1").Value = 'Datos actualizados a 19 de Mayo de 2020 a las 17:05'

# Estados Unidos (row 4) - updated totals
$ws.Range("B4").Value = 1554951
$ws.Range("C4").Value = 4657
$ws.Range("D4").Value = 359087
$ws.Range("E4").Value = 1103676
$ws.Range("G4").Value = 207
$ws.Range("H4").Value = 92188

# Alemania (row 11) - updated totals
$ws.Range("B11").Value = 177482
$ws.Range("C11").Value = 193
$ws.Range("E11").Value = 13637
$ws.Range("G11").Value = 22
$ws.Range("H11").Value = 8145

# India (row 14) - updated totals
$ws.Range("B14").Value = 103292
$ws.Range("C14").Value = 2964
$ws.Range("D14").Value = 40458
$ws.Range("E14").Value = 59655
$ws.Range("G14").Value = 23
$ws.Range("H14").Value = 3179

# Singapur (row 30) - updated active/recovered
$ws.Range("D30").Value = 10365
$ws.Range("E30").Value = 18407

# Rows 44-46: Republica Dominicana's new total (13223) overtakes both
# Filipinas and Egipto, so it moves up into row 44, pushing the other
# two down one row each (their own figures are unchanged).
$ws.Range("A44").Value = 'Republica Dominicana'
$ws.Range("B44").Value = 13223
$ws.Range("C44").Value = 498
$ws.Range("D44").Value = 6613
$ws.Range("E44").Value = 6169
$ws.Range("G44").Value = 7
$ws.Range("H44").Value = 441

$ws.Range("A45").Value = 'Filipinas'
$ws.Range("B45").Value = 12942
$ws.Range("C45").Value = 224
$ws.Range("D45").Value = 2843
$ws.Range("E45").Value = 9262
$ws.Range("G45").Value = 6
$ws.Range("H45").Value = 837

$ws.Range("A46").Value = 'Egipto'
$ws.Range("B46").Value = 12764
$ws.Range("D46").Value = 3440
$ws.Range("E46").Value = 8679
$ws.Range("H46").Value = 645

# Azerbaiyan (row 72) - updated totals
$ws.Range("B72").Value = 3518
$ws.Range("C72").Value = 131
$ws.Range("D72").Value = 2198
$ws.Range("E72").Value = 1279
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = 41

# Mozambique (row 158) - updated totals
$ws.Range("B158").Value = 146
$ws.Range("C158").Value = 1
$ws.Range("E158").Value = 102

# Rows 196-197: Nueva Caledonia and Belice swap rank (tied totals, but
# Nueva Caledonia's active-case count overtakes Belice's).
$ws.Range("A196").Value = 'Nueva Caledonia'
$ws.Range("D196").Value = 18
$ws.Range("H196").Value = 0

$ws.Range("A197").Value = 'Belice'
$ws.Range("D197").Value = 16
$ws.Range("H197").Value = 2

# Rows 209-211: Seychelles moves up ahead of Groenlandia and Montserrat
# (tied totals), shifting those two down one row each.
$ws.Range("A209").Value = 'Seychelles'

$ws.Range("A210").Value = 'Groenlandia'
$ws.Range("D210").Value = 11
$ws.Range("H210").Value = 0

$ws.Range("A211").Value = 'Montserrat'
$ws.Range("D211").Value = 10
$ws.Range("H211").Value = 1

# Rows 215-216: San Bartolome and Bonaire/San Eustaquio y Saba swap rank
# (tied totals).
$ws.Range("A215").Value = 'San Bartolome'
$ws.Range("A216").Value = 'Bonaire, San Eustaquio y Saba'
